$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are plain text in the source data (e.g. "69.422.72"
# using "." as a thousands separator), so force Text format before assigning
# to avoid Excel silently re-interpreting them as numbers (which would strip
# trailing zeros / switch to scientific notation for small values).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.406.99"
$ws.Range("E2").Value = "  -2.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.692.56"
$ws.Range("E3").Value = "  -3.02%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "691.94"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.10"
$ws.Range("E6").Value = "  -5.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.693.06"
$ws.Range("E7").Value = "  -3.00%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -4.78%  "

$ws.Range("E10").Value = "  -8.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.36"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("E12").Value = "  -5.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.34"
$ws.Range("E14").Value = "  -7.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.314.02"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.698.81"
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.450.60"
$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.15"
$ws.Range("E19").Value = "  -7.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.56"
$ws.Range("E20").Value = "  -8.06%  "

$ws.Range("E21").Value = "  -6.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.01"
$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.663"
$ws.Range("E23").Value = "  -7.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.85"
$ws.Range("E24").Value = "  -4.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.837.94"
$ws.Range("E25").Value = "  -3.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000129"
$ws.Range("E26").Value = "  -9.37%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("E28").Value = "  -5.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  -8.63%  "

$ws.Range("E30").Value = "  -11.09%  "

$ws.Range("E31").Value = "  -9.91%  "

$ws.Range("E32").Value = "  -7.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.06"
$ws.Range("E33").Value = "  -7.63%  "

$ws.Range("E34").Value = "  -5.63%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.19%  "

$ws.Range("E36").Value = "  -7.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.655.59"
$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.46"
$ws.Range("E38").Value = "  -7.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.32"
$ws.Range("E39").Value = "  +5.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.35"
$ws.Range("E40").Value = "  -1.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0926"
$ws.Range("E41").Value = "  -8.14%  "

$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.952"
$ws.Range("E44").Value = "  -6.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.28"
$ws.Range("E45").Value = "  -5.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "48.05"
$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.96"
$ws.Range("E47").Value = "  +2.77%  "

$ws.Range("E48").Value = "  -15.23%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.34"
$ws.Range("E49").Value = "  -2.67%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("E51").Value = "  -8.99%  "
